$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 (pushes "fossil_routes" and everything
# below it down by one row) and populate it with the new
# "chemical_recycling_pyrolysis" parameter, matching the row above it
# (chemical_recycling_gasification) in shape: column A = parameter name,
# column B = TRUE.
$ws.Rows(10).Insert()
$ws.Range("A10").Value = "chemical_recycling_pyrolysis"
$ws.Range("B10").Value = $true
